# Update the results table and chart title, and move the selection,
# matching the "updated results with new values" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 4 (intro_sort_many_equal): new benchmark results ---------------
$row4 = New-Object 'object[,]' 1,9
$row4[0,0] = 5480
$row4[0,1] = 9780
$row4[0,2] = 17600
$row4[0,3] = 33050
$row4[0,4] = 63930
$row4[0,5] = 123440
$row4[0,6] = 241910
$row4[0,7] = 477680
$row4[0,8] = 944370
$ws.Range("B4:J4").Value = $row4

# --- Row 5 (std_sort_many_equal): new benchmark results ------------------
$row5 = New-Object 'object[,]' 1,9
$row5[0,0] = 5620
$row5[0,1] = 10970
$row5[0,2] = 21310
$row5[0,3] = 42120
$row5[0,4] = 83760
$row5[0,5] = 164110
$row5[0,6] = 339750
$row5[0,7] = 691640
$row5[0,8] = 1415520
$ws.Range("B5:J5").Value = $row5

# --- Fix typo in the chart title -----------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$chart.ChartTitle.Text = "std::sort vs IntroSort bei vielen gleichen bzw. vielen unterschiedlichen Elementen"

# --- Move the active selection to T13 ------------------------------------
$ws.Range("T13").Select()
